$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the latest daily spot-price record; this automated update
# advances the date by one day and refreshes all hourly/summary values.

$ws.Range("A2").Value = 45886

$ws.Range("B2").Value = 118.58
$ws.Range("C2").Value = 113.12
$ws.Range("D2").Value = 106.27
$ws.Range("E2").Value = 104.59
$ws.Range("F2").Value = 104.66
$ws.Range("G2").Value = 105.01
$ws.Range("H2").Value = 104.34
$ws.Range("I2").Value = 96.98999999999999
$ws.Range("J2").Value = 101.12
$ws.Range("K2").Value = 27.2
$ws.Range("L2").Value = 9.699999999999999
$ws.Range("M2").Value = 2.01
$ws.Range("N2").Value = 0.65
$ws.Range("O2").Value = 0.65
$ws.Range("P2").Value = 0.65
$ws.Range("Q2").Value = 0.65
$ws.Range("R2").Value = 5.13
$ws.Range("S2").Value = 16.72
$ws.Range("T2").Value = 51.45
$ws.Range("U2").Value = 94.03
$ws.Range("V2").Value = 113.12
$ws.Range("W2").Value = 124.65
$ws.Range("X2").Value = 131.25
$ws.Range("Y2").Value = 117.67
$ws.Range("Z2").Value = 68.76000000000001

$ws.Range("AB2").Value = 121.67
$ws.Range("AC2").Value = "22h-24h"
$ws.Range("AD2").Value = 124.46
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 118.88
